# Refactor data parsing logic: append one newly-parsed telemetry row to each
# of the four per-station worksheets (ROW50-FE-LIFTER, ROW50-MID-LIFTER,
# ROW11-FE-LIFTER, ROW11-MID-LIFTER).

$wb = $excel.ActiveWorkbook
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# --- Sheet "ROW50-FE-LIFTER": append row 19 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$ws1.Cells.Item(19, 1).Value = 45733.62027259259
$ws1.Cells.Item(19, 1).NumberFormat = $dateFmt
$ws1.Cells.Item(19, 2).Value = "0x01,0x90"
$ws1.Cells.Item(19, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item(19, 4).Value = "0x01,0x82"
$ws1.Cells.Item(19, 5).Value = "0xe"
$ws1.Cells.Item(19, 6).Value = 400
$ws1.Cells.Item(19, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item(19, 8).Value = 386
$ws1.Cells.Item(19, 9).Value = 14

# --- Sheet "ROW50-MID-LIFTER": append row 21 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$ws2.Cells.Item(21, 1).Value = 45733.59885416667
$ws2.Cells.Item(21, 1).NumberFormat = $dateFmt
$ws2.Cells.Item(21, 2).Value = "0x01,0x90 "
$ws2.Cells.Item(21, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item(21, 4).Value = "0x01,0x86"
$ws2.Cells.Item(21, 5).Value = "0x19"
$ws2.Cells.Item(21, 6).Value = 400
# This value has more significant digits than a double can hold; prefix it
# with an apostrophe so Excel stores it as literal text instead of silently
# rounding it to scientific notation, then reset the style back to the
# sheet's default so only the value (not the formatting) changes.
$ws2.Cells.Item(21, 7).Value = "'568631262647113771663628"
$ws2.Cells.Item(21, 7).Style = "Normal"
$ws2.Cells.Item(21, 8).Value = 390
$ws2.Cells.Item(21, 9).Value = 25

# --- Sheet "ROW11-FE-LIFTER": append row 19 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$ws3.Cells.Item(19, 1).Value = 45733.64138497685
$ws3.Cells.Item(19, 1).NumberFormat = $dateFmt
$ws3.Cells.Item(19, 2).Value = "0x01,0x90"
$ws3.Cells.Item(19, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item(19, 4).Value = "0x01,0x82"
$ws3.Cells.Item(19, 5).Value = "0x14"
$ws3.Cells.Item(19, 6).Value = 400
$ws3.Cells.Item(19, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item(19, 8).Value = 386
$ws3.Cells.Item(19, 9).Value = 20

# --- Sheet "ROW11-MID-LIFTER": append row 19 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$ws4.Cells.Item(19, 1).Value = 45733.78845771991
$ws4.Cells.Item(19, 1).NumberFormat = $dateFmt
$ws4.Cells.Item(19, 2).Value = "0x01,0x90"
$ws4.Cells.Item(19, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item(19, 4).Value = "0x01,0x86"
$ws4.Cells.Item(19, 5).Value = "0x19"
$ws4.Cells.Item(19, 6).Value = 400
$ws4.Cells.Item(19, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item(19, 8).Value = 390
$ws4.Cells.Item(19, 9).Value = 25
